# edit.ps1 - Re-applies a PowerPoint "Design" (theme) change plus the
# resulting default-table-style change that PowerPoint performs on tables
# that still use the (now stale) theme-default table style.
#
# What happened in the authored edit:
#   1. The presentation's applied Design was switched from the
#      "Integral" (Red Violet) theme to the plain "Office Theme".
#      Concretely this rewrites the 12 theme colours used by the slide
#      master/theme (ppt/theme/theme1.xml -> a:clrScheme).
#   2. Because the table on slide 5 was still using the theme's
#      implicit/default table style, PowerPoint re-pointed it at the
#      new design's default table style GUID
#      ({5B0B2F53-3674-4A3A-991D-24C1C38B654B}).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Swap the Design's colour scheme from "Red Violet" to the standard
#    "Office" palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
# ---------------------------------------------------------------------
function RGBColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officePalette = @(
    (RGBColor 0x00 0x00 0x00),  # dk1
    (RGBColor 0xFF 0xFF 0xFF),  # lt1
    (RGBColor 0x44 0x54 0x6A),  # dk2
    (RGBColor 0xE7 0xE6 0xE6),  # lt2
    (RGBColor 0x5B 0x9B 0xD5),  # accent1
    (RGBColor 0xED 0x7D 0x31),  # accent2
    (RGBColor 0xA5 0xA5 0xA5),  # accent3
    (RGBColor 0xFF 0xC0 0x00),  # accent4
    (RGBColor 0x44 0x72 0xC4),  # accent5
    (RGBColor 0x70 0xAD 0x47),  # accent6
    (RGBColor 0x05 0x63 0xC1),  # hlink
    (RGBColor 0x95 0x4F 0x72)   # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officePalette[$i - 1]
}

# ---------------------------------------------------------------------
# 2) The table on slide 5 picks up the new design's default table
#    style.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{5B0B2F53-3674-4A3A-991D-24C1C38B654B}")
    }
}
